$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D9").Value = -7.561
$ws.Range("D13").Value = -8.220000000000001
$ws.Range("D16").Value = -8.747
$ws.Range("D18").Value = -8.548000000000002
$ws.Range("D20").Value = -7.43
$ws.Range("D26").Value = -8.122999999999999
$ws.Range("D27").Value = -8.815999999999999
$ws.Range("D29").Value = -7.369999999999999
$ws.Range("D35").Value = -7.873
$ws.Range("D36").Value = -7.784000000000001
$ws.Range("D45").Value = -7.558000000000002
$ws.Range("D55").Value = -8.456999999999999
$ws.Range("D57").Value = -8.315000000000001
$ws.Range("D69").Value = -7.160999999999999
$ws.Range("D76").Value = -7.806999999999999
$ws.Range("D78").Value = -7.878
$ws.Range("D82").Value = -8.418000000000001
$ws.Range("D83").Value = -8.354000000000001
$ws.Range("D93").Value = -7.452
$ws.Range("D97").Value = -8.273
